# edit.ps1 - apply the sys-log.docx changes described by the commit
# "fix command names, clarify startup state for prepackaged VM"

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Text-box size tweak (wp:extent / a:ext) -----------------------------
# The Word object model keeps the outer "wp:extent" and inner "a:ext"
# synced to the same value, so we match the (more visually significant)
# outer extent exactly.
$shp = $d.Shapes(1)
$shp.Width = 5468620 / 12700.0
$shp.Height = 481330 / 12700.0

# --- Paragraph text fixes --------------------------------------------------

Replace-Text "Boot your Linux system or VM, log in, and then open a terminal window and start the lab:" `
    "Boot your Linux system or VM.  If necessary, log in and then open a terminal window and cd to the labtainer/labtainer-student directory.  The pre-packaged Labtainer VM will start with such a terminal open for you.   Then start the lab:"

Replace-Text "cd labtainer/labtainer-student" ""

Replace-Text "start.py sys-log" "labtainer sys-log"

Replace-Text "It may help to stretch the resulting bash terminal window to the right to provide for more output space.  Note the terminal displays the paths to two files on your Linux host: " `
    "Note the terminal displays the paths to two files on your Linux host: "

Replace-Text "If you chose to edit the lab report on a different system, you are responsible for copying the completed report back to the displayed path on your Linux system before using “stop.py” to stop the lab for the last time." `
    "If you chose to edit the lab report on a different system, you are responsible for copying the completed report back to the displayed path on your Linux system before using “stoplab” to stop the lab for the last time."

# Merge the three split runs "/etc/" + "rsyslog.d/50-default" + ".conf"
# into a single run (occurs twice in the document).
Replace-Text "/etc/rsyslog.d/50-default.conf" "/etc/rsyslog.d/50-default.conf"

Replace-Text "stop.py sys-log" "stoplab sys-log"

Replace-Text "If you modified the lab report on a different system, you must copy that completed file into the directory path displayed when you started the lab, and you must do that before typing “./stop.py”.   When you stop the lab, the system will display a path to the zipped lab results on your Linux system.  " `
    "If you modified the lab report on a different system, you must copy that completed file into the directory path displayed when you started the lab, and you must do that before typing “stoplab”.   When you stop the lab, the system will display a path to the zipped lab results on your Linux system.  "

# --- Table cell margin ------------------------------------------------------
$t = $d.Tables(1)
$t.LeftPadding = 153 / 20.0

# --- New ListLabel109..ListLabel126 character styles -----------------------
function Add-ListLabelStyle($num, $asciiHAnsi, $cs) {
    $s = $d.Styles.Add("ListLabel " + $num, 2)
    if ($asciiHAnsi) {
        $s.Font.Name = $asciiHAnsi
    }
    $s.Font.NameBi = $cs
    $s.QuickStyle = $true
}

Add-ListLabelStyle 109 "Courier New" "Symbol"
Add-ListLabelStyle 110 $null "Courier New"
Add-ListLabelStyle 111 $null "Wingdings"
Add-ListLabelStyle 112 $null "Symbol"
Add-ListLabelStyle 113 $null "Courier New"
Add-ListLabelStyle 114 $null "Wingdings"
Add-ListLabelStyle 115 $null "Symbol"
Add-ListLabelStyle 116 $null "Courier New"
Add-ListLabelStyle 117 $null "Wingdings"
Add-ListLabelStyle 118 "Courier New" "Symbol"
Add-ListLabelStyle 119 $null "Courier New"
Add-ListLabelStyle 120 $null "Wingdings"
Add-ListLabelStyle 121 $null "Symbol"
Add-ListLabelStyle 122 $null "Courier New"
Add-ListLabelStyle 123 $null "Wingdings"
Add-ListLabelStyle 124 $null "Symbol"
Add-ListLabelStyle 125 $null "Courier New"
Add-ListLabelStyle 126 $null "Wingdings"

Write-Output "edit.ps1 complete"
